# Update "人数" (count) column F values on the "展览" and "全部类型" sheets
# to reflect the latest generated output, per commit "Update gh-pages to
# output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F2").Value = 650
$ws1.Range("F3").Value = 743
$ws1.Range("F4").Value = 955
$ws1.Range("F5").Value = 740
$ws1.Range("F6").Value = 846
$ws1.Range("F7").Value = 412
$ws1.Range("F8").Value = 621
$ws1.Range("F9").Value = 136
$ws1.Range("F10").Value = 1227
$ws1.Range("F11").Value = 651
$ws1.Range("F13").Value = 519
$ws1.Range("F16").Value = 612
$ws1.Range("F17").Value = 4
$ws1.Range("F18").Value = 366
$ws1.Range("F19").Value = 358
$ws1.Range("F21").Value = 558
$ws1.Range("F22").Value = 96
$ws1.Range("F25").Value = 813

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F4").Value = 650
$ws4.Range("F7").Value = 743
$ws4.Range("F8").Value = 955
$ws4.Range("F9").Value = 740
$ws4.Range("F10").Value = 846
$ws4.Range("F11").Value = 412
$ws4.Range("F12").Value = 621
$ws4.Range("F13").Value = 136
$ws4.Range("F14").Value = 1227
$ws4.Range("F15").Value = 651
$ws4.Range("F19").Value = 519
$ws4.Range("F23").Value = 613
$ws4.Range("F25").Value = 4
$ws4.Range("F26").Value = 366
$ws4.Range("F27").Value = 358
$ws4.Range("F31").Value = 558
$ws4.Range("F36").Value = 96
$ws4.Range("F39").Value = 813
